$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos data refresh: update Price (D) and Volume(1h) (E) columns,
# and the two swapped rows (35/36: WEMIXToken <-> Hedera) across B:E.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.227.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.198.78'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.93%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.85'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.15'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.509'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.81%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.466'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0767'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.97'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.82'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -11.24%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.542.93'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.20'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.87'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.204.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.707'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.171.47'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0866'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.69'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.59'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.23'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '223.92'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.98%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.79'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.49'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.30%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.76'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.56'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -10.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.77'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.23%  '

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.33'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.82%  '

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0690'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.109'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.22'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0950'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.61'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.63'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.59'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.895.23'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.07'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -8.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0258'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.64%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -10.44%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '71.43'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.408.71'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.85'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.49%  '
